$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark "Fix max righe/colonne" task as Done, with completion date ---
$ws.Range("C3").Value = "X"
$ws.Range("D8").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 44933

# --- Mark "0. creare pagina costruzione dataset" task as Done, with completion date ---
$ws.Range("C4").Value = "X"
$ws.Range("D8").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = 44933

# --- Add note about missing URL ---
$ws.Range("E4").Value = "Mancante: URL"

# --- Row heights (wrapped-text rows reflowed) ---
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 180
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 45
$ws.Rows.Item(23).RowHeight = 45
$ws.Rows.Item(25).RowHeight = 45
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 120

# --- Move selection to the newly-annotated cell ---
$ws.Range("E4").Select() | Out-Null
